$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

# --- Insert a new row 2 (above SYLVERSON) for account 005055865 / G3C / 200000 ---
$ws.Range("2:2").Insert()
$ws.Cells.Item(2, 1).Value = "'005055865"
$ws.Cells.Item(2, 2).Value = "G3C"
$ws.Cells.Item(2, 3).Value = 200000

# --- SYLVERSON's balance (now row 3) drops from 186307.41 to 186000 ---
$ws.Cells.Item(3, 3).Value = 186000

# --- Insert three new rows at 5:7 (above RODRIGO 004392159, now pushed to row 8) ---
$ws.Range("5:7").Insert()
$ws.Cells.Item(5, 1).Value = "'004212438"
$ws.Cells.Item(5, 2).Value = "KENIA"
$ws.Cells.Item(5, 3).Value = 21190.23

$ws.Cells.Item(6, 1).Value = "'004267976"
$ws.Cells.Item(6, 2).Value = "E3"
$ws.Cells.Item(6, 3).Value = 14647.49

$ws.Cells.Item(7, 1).Value = "'004452912"
$ws.Cells.Item(7, 2).Value = "BRUNO"
$ws.Cells.Item(7, 3).Value = 7151.71

# --- Insert a new row 10 (below BRUNO 004754056, above NATHALIA 005547702) ---
$ws.Range("10:10").Insert()
$ws.Cells.Item(10, 1).Value = "'005270025"
$ws.Cells.Item(10, 2).Value = "DENIZE"
$ws.Cells.Item(10, 3).Value = 100
